$d = $word.ActiveDocument

function Set-DateCell($table, $row, $col, $value) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $value

    $full = $table.Cell($row, $col).Range
    $full.Font.Name = "Cambria"
    $full.Font.Size = 10
    $full.Font.Bold = 0
    $full.Font.Italic = 0
    $full.Font.AllCaps = 0
    $full.Font.SmallCaps = 0
    $full.Font.StrikeThrough = 0
    $full.Font.DoubleStrikeThrough = 0
    $full.Font.Color = 0
    $full.Font.Underline = 0
    $full.Font.UnderlineColor = 0
    $full.Font.Subscript = 0
    $full.Font.Superscript = 0
    $full.ParagraphFormat.Alignment = 0
}

# --- Table edits: fill in the "Action Date" cells (column 3) for the
# "Analysis" (Germonda) and "Preparing presentation" rows ---

$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $taskText = $tbl.Cell($r, 1).Range.Text
    if ($taskText.StartsWith("Analysis")) {
        Set-DateCell $tbl $r 3 "01/11/16"
    }
    elseif ($taskText.StartsWith("Preparing presentation")) {
        Set-DateCell $tbl $r 3 "03/11/16"
    }
}

# --- Text edit: split "...you faced during..." into
# "...you face" + "D" (bold) + "d during..." (bold) ---
$d.Content.Find.Execute(
    "4.Describe the difficulties you faced during the collaboration (and solutions if appreciable)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "4.Describe the difficulties you faceDd during the collaboration (and solutions if appreciable)",
    2)
